$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B2" = 0.2330827067669173
    "C2" = 0.4962406015037594
    "J2" = 0.01879699248120301
    "O2" = 0.007518796992481203
    "P2" = 0.1691729323308271
    "S2" = 0.07518796992481203
    "B3" = 0.02158273381294964
    "C3" = 0.04316546762589928
    "J3" = 0.05035971223021583
    "P3" = 0.6546762589928058
    "S3" = 0.2302158273381295
    "J4" = 0.07692307692307693
    "P4" = 0.6153846153846154
    "S4" = 0.3076923076923077
    "B6" = 0.03225806451612903
    "D6" = 0.009216589861751152
    "F6" = 0.06451612903225806
    "J6" = 0.2811059907834101
    "O6" = 0.02764976958525346
    "Q6" = 0.1105990783410138
    "R6" = 0.07834101382488479
    "S6" = 0.3963133640552995
    "B7" = 0.09941520467836257
    "D7" = 0.01169590643274854
    "F7" = 0.08187134502923976
    "J7" = 0.1228070175438596
    "O7" = 0.02339181286549707
    "Q7" = 0.1637426900584795
    "R7" = 0.1169590643274854
    "S7" = 0.3801169590643275
    "B8" = 0.08733624454148471
    "D8" = 0.01091703056768559
    "F8" = 0.0851528384279476
    "J8" = 0.1244541484716157
    "O8" = 0.02838427947598253
    "Q8" = 0.1397379912663755
    "R8" = 0.08078602620087336
    "S8" = 0.4432314410480349
    "B9" = 0.06201550387596899
    "D9" = 0.0310077519379845
    "F9" = 0.07751937984496124
    "J9" = 0.1472868217054264
    "O9" = 0.0310077519379845
    "Q9" = 0.1395348837209302
    "R9" = 0.1007751937984496
    "S9" = 0.4108527131782946
    "B10" = 0.1073943661971831
    "D10" = 0.02288732394366197
    "F10" = 0.07130281690140845
    "J10" = 0.1091549295774648
    "O10" = 0.02200704225352113
    "Q10" = 0.1936619718309859
    "R10" = 0.07746478873239436
    "S10" = 0.3961267605633803
    "G11" = 0.1544401544401544
    "J11" = 0.1042471042471042
    "K11" = 0.2239382239382239
    "L11" = 0.5057915057915058
    "S11" = 0.01158301158301158
    "G12" = 0.8208955223880597
    "J12" = 0.1268656716417911
    "L12" = 0.02238805970149254
    "S12" = 0.02985074626865672
    "G13" = 0.7894736842105263
    "J13" = 0.1842105263157895
    "S13" = 0.02631578947368421
    "F15" = 0.01052631578947368
    "H15" = 0.1684210526315789
    "I15" = 0.03157894736842105
    "J15" = 0.3368421052631579
    "K15" = 0.07368421052631578
    "M15" = 0.005263157894736842
    "O15" = 0.03684210526315789
    "S15" = 0.3368421052631579
    "F16" = 0.00641025641025641
    "H16" = 0.2371794871794872
    "I16" = 0.03205128205128205
    "J16" = 0.3974358974358974
    "K16" = 0.1282051282051282
    "M16" = 0.01923076923076923
    "O16" = 0.04487179487179487
    "S16" = 0.1346153846153846
    "F17" = 0.02023121387283237
    "H17" = 0.1878612716763006
    "I17" = 0.07803468208092486
    "J17" = 0.4421965317919075
    "K17" = 0.07225433526011561
    "M17" = 0.01734104046242774
    "N17" = 0.002890173410404624
    "O17" = 0.05202312138728324
    "S17" = 0.1271676300578035
    "F18" = 0.02824858757062147
    "H18" = 0.1977401129943503
    "I18" = 0.05084745762711865
    "J18" = 0.384180790960452
    "K18" = 0.1016949152542373
    "M18" = 0.01129943502824859
    "O18" = 0.06214689265536723
    "S18" = 0.1638418079096045
    "F19" = 0.0145985401459854
    "H19" = 0.2327656123276561
    "I19" = 0.06731549067315491
    "J19" = 0.3763179237631792
    "K19" = 0.09975669099756691
    "M19" = 0.0218978102189781
    "O19" = 0.0575831305758313
    "S19" = 0.129764801297648
}

foreach ($key in $values.Keys) {
    $ws.Range($key).Value = $values[$key]
}
